# CUS15: actualizacion de servicios, scripts y archivos de cotizacion
# Adds the "Cotizacion N°" label and its value next to the document title,
# and updates the active sheet view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COTIZACION")

# --- New cells: C1 (label) and D1 (value) -------------------------------
# C1 should look like the other section labels (A3 "Datos del Cliente",
# A9 "Datos del Emisor (Proveedor)"): bold, 12pt, dark blue.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C1").Value = "Cotizacion N°"
$ws.Range("D1").Value = 21630

# --- Sheet view: scroll position and active selection -------------------
$window = $excel.ActiveWindow
$window.ScrollRow = 4
$window.ScrollColumn = 1
$ws.Range("D10").Select() | Out-Null
